$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.473.87'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').Value = '1.948.83'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.65'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.91'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.39%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.373'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0782'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.82%  '
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.14'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.23%  '
$ws.Range('D13').Value = '2.238.82'
$ws.Range('E13').Value = '  +0.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.824'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.39'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.40%  '
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('D17').Value = '1.954.26'
$ws.Range('E17').Value = '  +0.73%  '
$ws.Range('D18').Value = '36.346.10'
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.15'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('D20').Value = '0.0₃0845'
$ws.Range('E20').Value = '  -2.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '228.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.04'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.85%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('E24').Value = '  +2.89%  '
$ws.Range('E25').Value = '  +1.81%  '
$ws.Range('E26').Value = '  +7.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.11'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '159.57'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.19'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E30').Value = '  +18.36%  '
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.71'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.69%  '
$ws.Range('E33').Value = '  -1.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.41'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.37%  '
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.45'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.25'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.63%  '
$ws.Range('E38').Value = '  -1.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.42'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -13.07%  '
$ws.Range('E40').Value = '  +0.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0952'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.36%  '
$ws.Range('E42').Value = '  +0.90%  '
$ws.Range('E43').Value = '  -0.48%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.354.71'
$ws.Range('E44').Value = '  +0.93%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.62'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.19'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.81%  '
$ws.Range('E47').Value = '  -0.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.08'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.82'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.14'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.29%  '
$ws.Range('D51').Value = '2.128.96'
$ws.Range('E51').Value = '  +0.41%  '
